$wb = $excel.ActiveWorkbook

# --- "About" sheet updates ---
$about = $wb.Worksheets.Item("About")

# Update source citation block (B3:B7)
$about.Range("B3").Value = "Massachusetts Institute of Technology"
$about.Range("B4").Value = 2021
$about.Range("B5").Value = "Re-examining rates of lithium-ion battery technology improvement and cost decline"
$about.Range("B6").Value = "https://pubs.rsc.org/en/content/articlepdf/2021/ee/d0ee02681f?page=search"
$about.Range("B7").Value = "Abstract"

# Clear the old note that referenced the now-removed chart (keep cell/style)
$about.Range("C8").Value = ""

# Add the new note about the averaged learning rate
$about.Range("A9").Value = "Note: We take the average of learning rates quoted in the Abstract (20%-27%)"

# Remove the embedded chart picture (no longer referenced)
if ($about.Shapes.Count -gt 0) {
    foreach ($shp in @($about.Shapes)) {
        $shp.Delete()
    }
}

# --- "PDiBCpDoC" sheet updates ---
$data = $wb.Worksheets.Item("PDiBCpDoC")

# Replace the hard-coded value with a formula averaging the two quoted rates
$data.Range("B2").Formula = "=AVERAGE(0.2,0.27)"

# Update saved selections to match the latest editing session
[void]$data.Range("I4").Select()
[void]$about.Range("A10").Select()
